# The presentation ships two theme parts:
#   ppt/theme/theme1.xml  - "Office Theme" colours, used by the notes master
#   ppt/theme/theme2.xml  - "Integral" colours, used by the slide master /
#                            the presentation's single Design ("Integral")
#
# The target revision swaps the two themes' contents: theme2.xml (the
# slide master's theme, which is what Designs/SlideMaster/Slides expose
# through the PowerPoint object model) becomes the "Office Theme" colour
# set. Apply that by rewriting every slot of the shared ThemeColorScheme
# that backs the presentation's Design/SlideMaster/Slides.

function ColorToCom([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target ("Office Theme") colours in ThemeColorScheme index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ColorToCom $officeThemeColors[$i - 1]
}
